# Weekly update: a new price-survey row is inserted for
# "Feria Lagunitas de Puerto Montt - Zanahoria" at sheet row 348,
# pushing the existing rows 348-371 down to 349-372 (dimension grows
# from A1:R371 to A1:R372).
#
# The inserted row duplicates the data currently sitting in row 348
# (same Volumen/Precio/Unidad/Precio-Kg figures) but is dated later
# and sourced from a different Origen ("Chillán").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 348:371 down one slot, leaving a blank row 348.
$ws.Rows.Item(348).Insert()

# Seed the new row 348 with the same record as the row right below it
# (which now holds what used to be row 348), then tweak the two cells
# that actually differ for the new entry.
$srcRow = $ws.Range("A349:R349")
$newRow = $ws.Range("A348:R348")
$newRow.Value2 = $srcRow.Value2

$ws.Range("D348").Value2 = 44714
$ws.Range("O348").Value = "Chillán"
